$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NestTable")

# Row 1 was previously empty (data started at row 2); just populate A1,
# no row-shifting is involved.
$ws.Range("A1").Value = "TableRowBase=FacilityTable"

$ws.Range("G17").Select()
